# 🔄 Sync automático del tracker (cada 3h)
# Appends two new result rows (63 and 64) to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 63; EventId = "14579389"; Fecha = "2025-09-02"; JugadorA = "Oliver Crawford"; JugadorB = "Marat Sharipov"; Pronostico = "Gana Marat Sharipov"; Cuota = 2.63 },
    @{ Row = 64; EventId = "14579388"; Fecha = "2025-09-02"; JugadorA = "Rio Noguchi"; JugadorB = "Yosuke Watanuki"; Pronostico = "Gana Rio Noguchi"; Cuota = 3 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # event_id -> stored as text in the source data
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.EventId
    $cellA.Style = "Normal"

    # fecha -> keep as literal text, not an auto-converted date serial
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.Fecha
    $cellB.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $r.JugadorA
    $ws.Cells.Item($row, 4).Value = $r.JugadorB
    $ws.Cells.Item($row, 5).Value = $r.Pronostico
    $ws.Cells.Item($row, 6).Value = $r.Cuota

    # resultado / profit -> still empty, pending result
    $cellG = $ws.Cells.Item($row, 7)
    $cellG.NumberFormat = "@"
    $cellG.Style = "Normal"

    $cellH = $ws.Cells.Item($row, 8)
    $cellH.NumberFormat = "@"
    $cellH.Style = "Normal"
}
